# Injuries_Master_Clubs.xlsx refresh (2025-11-15 -> 2025-11-16 scrape)
#
# snapshot:   two players ("СИБ" Гордеев Фёдор, "ШДР" Саттер Райли) came off
#             the injury list, so their rows are removed and every remaining
#             row's scraped_at (col K) is bumped to the new scrape run.
# returned:   now lists those two players as RETURN events for 2025-11-16.
# new_injured: no new injuries in this run, so only the header row remains.

# NOTE: positional parameters only -- named-parameter calls (-Cell ... -Value
# ...) into this runtime's PowerShell host don't propagate COM mutations
# back to the caller's Range object, so the function must be called as
# `Set-TextCell <cell> <value>`.
function Set-TextCell($Cell, $Value) {
    # Force the cell to keep a literal string type (not auto-coerced to a
    # number/date by COM), then drop back to the default "Normal" style so
    # no stray NumberFormat/style survives on the cell.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "snapshot" sheet: drop the two recovered players, refresh timestamps
# ---------------------------------------------------------------------
$snapshot = $wb.Worksheets.Item("snapshot")

# Row 20 = СИБ / Гордеев Фёдор, Row 44 (becomes 43 after first delete) =
# ШДР / Саттер Райли. Deleting shifts everything below up automatically.
$snapshot.Rows.Item(20).Delete()
$snapshot.Rows.Item(43).Delete()

$timestamps = @(
    @{Row=2; Value="2025-11-16T03:01:41.167735+00:00"},
    @{Row=3; Value="2025-11-16T03:01:41.167773+00:00"},
    @{Row=4; Value="2025-11-16T03:01:41.167795+00:00"},
    @{Row=5; Value="2025-11-16T03:01:43.861639+00:00"},
    @{Row=6; Value="2025-11-16T03:01:43.861671+00:00"},
    @{Row=7; Value="2025-11-16T03:01:46.213765+00:00"},
    @{Row=8; Value="2025-11-16T03:01:48.917073+00:00"},
    @{Row=9; Value="2025-11-16T03:01:51.677663+00:00"},
    @{Row=10; Value="2025-11-16T03:01:51.677692+00:00"},
    @{Row=11; Value="2025-11-16T03:01:56.695444+00:00"},
    @{Row=12; Value="2025-11-16T03:01:59.038251+00:00"},
    @{Row=13; Value="2025-11-16T03:02:01.764048+00:00"},
    @{Row=14; Value="2025-11-16T03:02:01.764080+00:00"},
    @{Row=15; Value="2025-11-16T03:02:01.764099+00:00"},
    @{Row=16; Value="2025-11-16T03:02:04.581142+00:00"},
    @{Row=17; Value="2025-11-16T03:02:07.280322+00:00"},
    @{Row=18; Value="2025-11-16T03:02:07.280357+00:00"},
    @{Row=19; Value="2025-11-16T03:02:09.951485+00:00"},
    @{Row=20; Value="2025-11-16T03:02:12.322779+00:00"},
    @{Row=21; Value="2025-11-16T03:02:12.322812+00:00"},
    @{Row=22; Value="2025-11-16T03:02:12.322832+00:00"},
    @{Row=23; Value="2025-11-16T03:02:12.322849+00:00"},
    @{Row=24; Value="2025-11-16T03:02:15.147779+00:00"},
    @{Row=25; Value="2025-11-16T03:02:15.147808+00:00"},
    @{Row=26; Value="2025-11-16T03:02:17.921853+00:00"},
    @{Row=27; Value="2025-11-16T03:02:17.921884+00:00"},
    @{Row=28; Value="2025-11-16T03:02:17.921905+00:00"},
    @{Row=29; Value="2025-11-16T03:02:20.597059+00:00"},
    @{Row=30; Value="2025-11-16T03:02:20.597088+00:00"},
    @{Row=31; Value="2025-11-16T03:02:23.359413+00:00"},
    @{Row=32; Value="2025-11-16T03:02:23.359442+00:00"},
    @{Row=33; Value="2025-11-16T03:02:23.359460+00:00"},
    @{Row=34; Value="2025-11-16T03:02:23.359476+00:00"},
    @{Row=35; Value="2025-11-16T03:02:23.359490+00:00"},
    @{Row=36; Value="2025-11-16T03:02:26.167578+00:00"},
    @{Row=37; Value="2025-11-16T03:02:26.167607+00:00"},
    @{Row=38; Value="2025-11-16T03:02:31.264393+00:00"},
    @{Row=39; Value="2025-11-16T03:02:31.264423+00:00"},
    @{Row=40; Value="2025-11-16T03:02:31.264440+00:00"},
    @{Row=41; Value="2025-11-16T03:02:31.264455+00:00"},
    @{Row=42; Value="2025-11-16T03:02:33.640826+00:00"}
)

foreach ($item in $timestamps) {
    Set-TextCell $snapshot.Cells.Item($item.Row, 11) $item.Value
}

# ---------------------------------------------------------------------
# 2) "returned" sheet: replace the 3 old RETURN rows with the 2 new ones
# ---------------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")

# Drop the extra old data row (3 rows -> 2 rows); rows 2 & 3 get overwritten
# below with the new content.
$returned.Rows.Item(4).Delete()

$returnedRows = @(
    @{Row=2; Values=@("СИБ", "Сибирь", "Гордеев Фёдор", "1369_СИБ_гордеевфедор", "RETURN", "2025-11-16T11:02:34.145192+08:00", "2025-11-16")},
    @{Row=3; Values=@("ШДР", "Драконы", "Саттер Райли", "1369_ШДР_саттеррайли", "RETURN", "2025-11-16T11:02:34.145192+08:00", "2025-11-16")}
)

foreach ($rowSpec in $returnedRows) {
    for ($col = 1; $col -le 7; $col++) {
        Set-TextCell $returned.Cells.Item($rowSpec.Row, $col) $rowSpec.Values[$col - 1]
    }
}

# ---------------------------------------------------------------------
# 3) "new_injured" sheet: no new injuries this run -> header row only
# ---------------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
$newInjured.Rows.Item(2).Delete()
